$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Förändrad" (Changed) column C for every existing data
#        row (2..339) from 2023-09-11 (45180) to 2023-09-12 (45181). ---
$ws.Range("C2:C339").Value = 45181

# --- 2. Row 339 picks up an explicit row-height stamp (matches the rest
#        of the data rows, which already carry ht="15" customHeight="1"). ---
$ws.Rows.Item(339).RowHeight = 15

# --- 3. Append three new cleared-logging notifications as rows 340-342. ---
$newRows = @(
    @{ Row = 340; A = "A 42324-2023"; G = 2.2 },
    @{ Row = 341; A = "A 42328-2023"; G = 1.3 },
    @{ Row = 342; A = "A 42326-2023"; G = 1 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A                 # A: Beteckning
    $ws.Cells.Item($row, 2).Value = 45180                # B: Datum
    $ws.Cells.Item($row, 3).Value = 45181                # C: Förändrad
    $ws.Range("B$row`:C$row").NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 4).Value = "ÖREBRO LÄN"          # D: Län
    $ws.Cells.Item($row, 5).Value = "LAXÅ"                # E: Kommun
    $ws.Cells.Item($row, 6).Value = "Sveaskog"            # F: Markägare
    $ws.Cells.Item($row, 7).Value = $r.G                  # G: Area (ha)

    # H..Q: species / conservation-status counters, all zero for these rows
    $ws.Range("H$row`:Q$row").Value = 0

    # R: Artnamn column keeps the wrap-text style even though it's empty
    $ws.Range("R$row").WrapText = $true
}

# Rows 340 and 341 also get the explicit row-height stamp; row 342 (the
# new last row) stays without it, mirroring the previous "last row has no
# explicit height" pattern.
$ws.Rows.Item(340).RowHeight = 15
$ws.Rows.Item(341).RowHeight = 15
